# StoryAndTask.xlsx - "code quality and improvement" pass.
#
# Adds the missing Sprint-6 sub-tasks (S6T1 gets real text instead of a
# stub, plus brand-new S6T2-S6T6 rows), a new S4T4 task, fixes the S7T5
# wording, and adds a new S10T1 "success message" task. The shared-string
# table / cell indices re-pack themselves on save, so we only need to set
# the cell values that actually change or are brand new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 6: S6T1 was a stub ("S6T1: ") - give it real content. ---
$ws.Range("B45").Value = "S6T1: As a developer I need to create portfolio Button in UI(Android)"

# --- Sprint 4: new task S4T4, inserted right after S4T3 (row 32) ---
$ws.Range("B33").Value = "S4T4: As a developer I need to add finctinily so that succesfull message will appear for user."

# --- Sprint 6 (cont'd): five brand-new sub-tasks S6T2..S6T6 under S6T1. ---
$ws.Range("B46").Value = "S6T2: As a developer I need to set on click listener in android."
$ws.Range("B47").Value = "S6T3: As a developer I need to create new screen to show all details of order submission."
$ws.Range("B48").Value = "S6T4: As a developer I need to create search bar so that user can enter account id to see particular client portfolio."
$ws.Range("B49").Value = "S6T5: As a developer I need to create functinality so that information retirve from database and show on UI with REST"
$ws.Range("B50").Value = "S6T6: AS a developer I need to create"

# --- Sprint 7: fix typo / reword S7T5 ---
$ws.Range("B60").Value = "S7T5: As a developer I need to add some functionality  on button to disbale and enble accordingly."

# --- Sprint 10: new task, a second S10T1 entry about the logout success message ---
$ws.Range("B87").Value = "S10T1: As a developer I need to add a message for successfully logout."

# Leave the freshly added row selected, matching the author's final cursor position.
$ws.Range("B87").Select()
